# Update the "Instructions" sheet:
#  - bump the version string
#  - split the old combined "Add antibodies / don't edit other sheets" line
#    into three separate instruction lines (one new row added)
#  - this shifts every row below the instructions block down by two rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected (no password), so it must be unprotected before
# any cell contents can be changed.
$ws.Unprotect()

# Make room for the two extra instruction lines (new rows 5 and 6); this
# pushes the old row 5 (blank) onward down to row 7, old row 6 ("Antibody
# name") down to row 8, etc.
$ws.Range("A5:A6").EntireRow.Insert()

# Bump the version number.
$ws.Range("A2").Value = "Version 1.2.2"

# Re-split the instructions text across three rows.
$ws.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet."
$ws.Range("A5").Value = "Do not change the headers of the 'Antibodies' sheet."
$ws.Range("A6").Value = "Do not edit the other sheets."

# Restore sheet protection.
$ws.Protect()
